$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.86075496673584
$ws.Range("B1").Value = 1.884145617485046
$ws.Range("C1").Value = 1.673977851867676
$ws.Range("D1").Value = 1.745153665542603
$ws.Range("E1").Value = 1.675243377685547
